$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1303
$ws.Range("C2").Value = 1887

$ws.Range("C3").Value = 66

$ws.Range("B4").Value = 269
$ws.Range("C4").Value = 354

$ws.Range("B5").Value = 107
$ws.Range("C5").Value = 153
